$wb = $excel.ActiveWorkbook

# Mapping of row number -> new value for column F ("想去人数")
$updates = @{
    3  = 1003
    4  = 238
    5  = 1366
    6  = 8406
    7  = 51
    10 = 233
    11 = 142
    12 = 3374
    14 = 336
    16 = 897
    17 = 137
    18 = 1088
    20 = 144
    21 = 1962
}

# Both the "展览" and "全部类型" sheets contain the same data and both
# need the column F values updated.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
